$d = $word.ActiveDocument
$t = $d.Tables.Item(1)
$values = @(
  "32+40=72",
  "43-33=10",
  "13+24=37",
  "62+20=82",
  "12-1=11",
  "46-9=37",
  "97-54=43",
  "84+7=91",
  "79-46=33",
  "4+57=61",
  "35-0=35",
  "8+22=30",
  "67-62=5",
  "85-53=32",
  "72-44=28",
  "76+9=85",
  "48+35=83",
  "77-53=24",
  "63-53=10",
  "35-23=12",
  "37+29=66",
  "56-2=54",
  "38+59=97",
  "97-66=31",
  "40+35=75",
  "65-62=3",
  "67-42=25",
  "12+36=48",
  "34+62=96",
  "30+60=90",
  "79-22=57",
  "50-25=25",
  "24+20=44",
  "54-19=35",
  "46+25=71",
  "33+1=34",
  "53-37=16",
  "78-24=54",
  "72-63=9",
  "47+22=69",
  "44+3=47",
  "99-2=97",
  "25+68=93",
  "18+3=21",
  "81-41=40",
  "38-2=36",
  "16+31=47",
  "89-78=11",
  "99-90=9",
  "89-31=58",
  "46+14=60",
  "8-2=6",
  "20+14=34",
  "58-4=54",
  "56+24=80",
  "85+6=91",
  "18-13=5",
  "23+29=52",
  "88-7=81",
  "40+38=78",
  "84-33=51",
  "59-5=54",
  "5+48=53",
  "27-26=1",
  "55+22=77",
  "20+49=69",
  "35+45=80",
  "71+6=77",
  "69-50=19",
  "30-4=26",
  "26+32=58",
  "26+46=72",
  "7+8=15",
  "9+29=38",
  "56+36=92",
  "4+47=51",
  "25+43=68",
  "18+57=75",
  "93-84=9",
  "13-3=10",
  "82-37=45",
  "38-9=29",
  "9-9=0",
  "99-66=33",
  "33+18=51",
  "50+41=91",
  "10-8=2",
  "23+32=55",
  "94-1=93",
  "71+28=99",
  "54-6=48",
  "61+28=89",
  "67+12=79",
  "20+73=93",
  "34+24=58",
  "95-75=20",
  "26+59=85",
  "5+31=36",
  "2+35=37",
  "76-20=56"
)

$expected = $values.Count
$total = $t.Rows.Count * $t.Columns.Count
if ($total -ne $expected) {
  Write-Host "WARNING: table has" $total "cells, expected" $expected
}

$idx = 0
for ($r = 1; $r -le $t.Rows.Count; $r++) {
  for ($c = 1; $c -le $t.Columns.Count; $c++) {
    if ($idx -lt $values.Count) {
      $cell = $t.Cell($r, $c)
      $cell.Range.Text = $values[$idx]
    }
    $idx = $idx + 1
  }
}
Write-Host "Replaced" $idx "cells"
